# Auto update Excel log
# Appends newly captured sensor readings to the PIR, Humidity, and Temperature sheets

$wb = $excel.ActiveWorkbook

$pirRows = @(
    ,@(80, "2026-01-28", "16:17:05", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@(81, "2026-01-28", "16:17:06", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@(82, "2026-01-28", "16:17:10", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@(83, "2026-01-28", "16:17:15", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@(84, "2026-01-28", "16:17:20", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@(85, "2026-01-28", "16:17:25", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@(86, "2026-01-28", "16:17:30", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@(87, "2026-01-28", "16:17:35", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@(88, "2026-01-28", "16:17:40", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@(89, "2026-01-28", "16:17:45", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@(90, "2026-01-28", "16:17:50", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@(91, "2026-01-28", "16:17:55", "16:00", "Bathroom", "No Motion", "Inactive")
    ,@(92, "2026-01-28", "16:18:00", "16:00", "Bathroom", "No Motion", "Inactive")
)

$humRows = @(
    ,@(79, "2026-01-28", "16:17:05", "16:00", "Bathroom", "88.3%", "Active")
    ,@(80, "2026-01-28", "16:17:06", "16:00", "Bathroom", "87.4%", "Active")
    ,@(81, "2026-01-28", "16:17:07", "16:00", "Bathroom", "88.3%", "Active")
    ,@(82, "2026-01-28", "16:17:11", "16:00", "Bathroom", "88.3%", "Active")
    ,@(83, "2026-01-28", "16:17:15", "16:00", "Bathroom", "87.5%", "Active")
    ,@(84, "2026-01-28", "16:17:19", "16:00", "Bathroom", "88.4%", "Active")
    ,@(85, "2026-01-28", "16:17:23", "16:00", "Bathroom", "87.5%", "Active")
    ,@(86, "2026-01-28", "16:17:27", "16:00", "Bathroom", "88.3%", "Active")
    ,@(87, "2026-01-28", "16:17:31", "16:00", "Bathroom", "88.4%", "Active")
    ,@(88, "2026-01-28", "16:17:35", "16:00", "Bathroom", "87.5%", "Active")
    ,@(89, "2026-01-28", "16:17:39", "16:00", "Bathroom", "88.3%", "Active")
    ,@(90, "2026-01-28", "16:17:47", "16:00", "Bathroom", "88.4%", "Active")
    ,@(91, "2026-01-28", "16:17:51", "16:00", "Bathroom", "88.4%", "Active")
    ,@(92, "2026-01-28", "16:17:55", "16:00", "Bathroom", "87.5%", "Active")
    ,@(93, "2026-01-28", "16:17:59", "16:00", "Bathroom", "86.9%", "Active")
)

$tempRows = @(
    ,@(79, "2026-01-28", "16:17:05", "16:00", "Bathroom", "22.8C", "Active")
    ,@(80, "2026-01-28", "16:17:06", "16:00", "Bathroom", "22.8C", "Active")
    ,@(81, "2026-01-28", "16:17:07", "16:00", "Bathroom", "22.7C", "Active")
    ,@(82, "2026-01-28", "16:17:11", "16:00", "Bathroom", "22.8C", "Active")
    ,@(83, "2026-01-28", "16:17:15", "16:00", "Bathroom", "22.8C", "Active")
    ,@(84, "2026-01-28", "16:17:19", "16:00", "Bathroom", "22.8C", "Active")
    ,@(85, "2026-01-28", "16:17:23", "16:00", "Bathroom", "22.8C", "Active")
    ,@(86, "2026-01-28", "16:17:27", "16:00", "Bathroom", "22.7C", "Active")
    ,@(87, "2026-01-28", "16:17:31", "16:00", "Bathroom", "22.8C", "Active")
    ,@(88, "2026-01-28", "16:17:35", "16:00", "Bathroom", "22.8C", "Active")
    ,@(89, "2026-01-28", "16:17:39", "16:00", "Bathroom", "22.7C", "Active")
    ,@(90, "2026-01-28", "16:17:47", "16:00", "Bathroom", "22.8C", "Active")
    ,@(91, "2026-01-28", "16:17:51", "16:00", "Bathroom", "22.8C", "Active")
    ,@(92, "2026-01-28", "16:17:55", "16:00", "Bathroom", "22.8C", "Active")
    ,@(93, "2026-01-28", "16:18:00", "16:00", "Bathroom", "22.7C", "Active")
)


function Write-SensorRows {
    param($ws, $rows)
    foreach ($row in $rows) {
        $r = $row[0]
        $dateVal = $row[1]
        $timeVal = $row[2]
        $hourVal = $row[3]
        $locVal = $row[4]
        $valueVal = $row[5]
        $statusVal = $row[6]

        $ws.Cells.Item($r, 1).Value = "'" + $dateVal
        $ws.Cells.Item($r, 2).Value = $timeVal
        $ws.Cells.Item($r, 3).Value = $hourVal
        $ws.Cells.Item($r, 4).Value = $locVal

        if ($valueVal -match '^[0-9.]+%$') {
            $ws.Cells.Item($r, 5).Value = "'" + $valueVal
        } else {
            $ws.Cells.Item($r, 5).Value = $valueVal
        }

        $ws.Cells.Item($r, 6).Value = $statusVal
    }
}

$pirWs = $wb.Worksheets.Item("PIR")
Write-SensorRows $pirWs $pirRows

$humWs = $wb.Worksheets.Item("Humidity")
Write-SensorRows $humWs $humRows

$tempWs = $wb.Worksheets.Item("Temperature")
Write-SensorRows $tempWs $tempRows

Write-Host "Rows appended: PIR->92, Humidity->93, Temperature->93"

